$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Veda")

# Row 11 becomes "windoff": clear G11/H11, set I11/J11 to new values
$ws.Range("F11").Value = "windoff"
$ws.Range("G11").Value = $null
$ws.Range("H11").Value = $null
$ws.Range("I11").Value = 0.12393998695368556
$ws.Range("J11").Value = 0.40029835465150038

# New row 12 "windon" carries the old wind values (G/H/I unchanged, J updated)
$ws.Range("F12").Value = "windon"
$ws.Range("G12").Value = 0.13646444879321595
$ws.Range("H12").Value = 0.24246061460306279
$ws.Range("I12").Value = 0.13634993213800239
$ws.Range("J12").Value = 0.21985176055114744

# Match number formatting/style of row 11 (style s="3") for the new row 12 cells
$ws.Range("G12:J12").NumberFormat = $ws.Range("G11:J11").NumberFormat

# The historical data sheet only tracked "wind" overall; now it is relabeled
# as "windon" to sit alongside the new "windoff" split introduced above.
$wsHist = $wb.Worksheets.Item("historical_data_long")
$wsHist.Range("A1:A698").Replace("wind", "windon", 1, 1, $false) | Out-Null

$wb.Application.Calculate()
